# SOFT153 CW.docx edit script
#
# Changes applied (per the provided unified diff):
#  1. Remove the `_GoBack` bookmark from the "Comparison" heading paragraph.
#  2. Split the run containing "O(n log n)" into several runs (adds a
#     <w:proofErr> pair around "O(" and a <w:lastRenderedPageBreak/> before
#     "log n)") -- a purely cosmetic run-split with no visible text change.
#  3. Append a new sentence after the existing final sentence:
#     "After some research it appears that my insertion sort is seemingly
#      more efficient than the average insertion sort which is O(n2)."
#     (the trailing 2 is a superscript run).
#  4. Re-insert the `_GoBack` bookmark at the very end of that same
#     paragraph (after the freshly appended text).
#
# We drive this with Range.InsertXML so that we get exact control over the
# resulting run/bookmark structure (Find/Replace collapses runs, which we
# don't want here).

$d = $word.ActiveDocument

function Get-ParagraphIndexByText($doc, [string]$needle) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $p = $doc.Paragraphs($i)
        if ($p.Range.Text.Contains($needle)) {
            return $i
        }
    }
    return -1
}

$wordOpenXmlHeader = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>'
$wordOpenXmlFooter = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

# ---------------------------------------------------------------------
# 1) Strip the _GoBack bookmark from the "Comparison" paragraph. Plain
#    text/XML replacement leaves a point-bookmark like _GoBack in place
#    (it's anchored to a position, not "inside" the replaced text), so it
#    has to be removed explicitly via the Bookmarks collection.
# ---------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ---------------------------------------------------------------------
# 2)-4) Rebuild the final paragraph: split the O(n log n) run, append the
#    new insertion-sort sentence, and move the _GoBack bookmark to the end.
# ---------------------------------------------------------------------
$finalIdx = Get-ParagraphIndexByText $d "more closely represents O(n log n)"
$finalPara = $d.Paragraphs($finalIdx)
$finalRange = $d.Range($finalPara.Range.Start, $finalPara.Range.End - 1)

$finalInner = '<w:p>' +
  '<w:r><w:t>As you can see when both insertion sort and quicksort are plotted on the same graph</w:t></w:r>' +
  '<w:r w:rsidR="007232F2"><w:t xml:space="preserve">, insertion follows the O(n) </w:t></w:r>' +
  '<w:r w:rsidR="002D248A"><w:t>time complexity pattern</w:t></w:r>' +
  '<w:r w:rsidR="007232F2"><w:t xml:space="preserve"> and quicksort clearly shows it more closely represents </w:t></w:r>' +
  '<w:proofErr w:type="gramStart"/>' +
  '<w:r><w:t>O(</w:t></w:r>' +
  '<w:proofErr w:type="gramEnd"/>' +
  '<w:r><w:t xml:space="preserve">n </w:t></w:r>' +
  '<w:r><w:lastRenderedPageBreak/><w:t>log n)</w:t></w:r>' +
  '<w:r w:rsidR="00AB5E4B"><w:t>.</w:t></w:r>' +
  '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
  '<w:r><w:t xml:space="preserve">After some research it appears that my insertion sort is seemingly more efficient than the average </w:t></w:r>' +
  '<w:r><w:t>insertion sort which is O(n</w:t></w:r>' +
  '<w:r><w:rPr><w:vertAlign w:val="superscript"/></w:rPr><w:t>2</w:t></w:r>' +
  '<w:r><w:t>)</w:t></w:r>' +
  '<w:r><w:t>.</w:t></w:r>' +
  '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' +
  '</w:p>'
$finalXml = $wordOpenXmlHeader + $finalInner + $wordOpenXmlFooter
$finalRange.InsertXML($finalXml)
